$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New explanatory text above the input table (row 4) ---
$ws.Range("A4").Value = "Important: values filled in the yellow cells should be in the same units as the units used internally by oXs (e.g. milliVolt or milliAmpere)"

# --- Update the measured sample values (rows 5 and 6) ---
$ws.Range("B5").Value = 110
$ws.Range("C5").Value = 18000
$ws.Range("D5").Value = 40000

$ws.Range("B6").Value = 1670
$ws.Range("C6").Value = 2150
$ws.Range("D6").Value = 2700

# --- New rows 17-19: quick "sanity check" block ---
# Note: cells are written in this order (A18 before A17) so that the
# shared-string table ends up with the same index assignment as the
# target workbook (10=Important, 11=transmit, 12=so when, 13=mv, 14=Please note).
$ws.Range("A18").Value = "oXs will transmit to the handset =>"

$ws.Range("A17").Value = "so when oXs will measure e.g. =>"
$ws.Range("B17").Value = 2700
$ws.Range("B17").Interior.Color = $ws.Range("B5").Interior.Color
$ws.Range("C17").Value = "mv"

$ws.Range("B18").Formula = "=B17*B13-B14"

$ws.Range("A19").Value = "Please note that the units on the handset can be different (e.g. A instead of mA)"

$ws.Range("A19").Select() | Out-Null

# --- Move / resize the embedded chart to its new anchor position ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 443.0565618848425
$co.Top = 168
$co.Width = 469.75
$co.Height = 221.25
